$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
# Copy formatting from H1 (the last existing header cell) so the new
# header cells match the existing bold/centered/bordered header style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-9, columns I and J
$values = @(
    @(2, 7, 8),
    @(3, 6, 6),
    @(4, 7, 8),
    @(5, 1, 5),
    @(6, 1, 6),
    @(7, 1, 4),
    @(8, 1, 4),
    @(9, 1, 2)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
